$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "BUZZER / on-board buzzer" annotation from row 53 (Digital pin 40)
# up to row 51 (Digital pin 38).
$ws.Range("D51").Value = "BUZZER"
$ws.Range("E51").Value = "on-board buzzer"

# Clear the old BRAKE_ON / on-board relay annotation that used to sit on
# row 52 (Digital pin 41) - it is being relocated to row 71 below.
$ws.Range("D52").ClearContents()
$ws.Range("E52").ClearContents()

# Clear the BUZZER / on-board buzzer annotation that used to sit on row 53
# (Digital pin 40) now that it has moved up to row 51.
$ws.Range("D53").ClearContents()
$ws.Range("E53").ClearContents()

# Add the BRAKE_ON annotation (relay wording tweaked to "on board relay")
# to row 71 (Digital pin 39).
$ws.Range("D71").Value = "BRAKE_ON"
$ws.Range("E71").Value = "on board relay"

# Update the active selection to match where the author ended up editing.
[void]$ws.Range("D72").Select()
